$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update existing numeric values (higher precision re-measurements) ---
$ws.Range("F2").Value2 = 0.97199999999999998
$ws.Range("F3").Value2 = 0.57799999999999996
$ws.Range("F4").Value2 = 0.72499999999999998

# --- Row 5 ("Nicht simultan" / "Clean"): fill in new statistics ---
$ws.Range("C5").Value2 = 1
$ws.Range("E5").Value2 = 1

# D5 / F5 need to hold the *text* values "0.81" / "0.977" (numbers-as-text,
# matching how similar cells elsewhere in the sheet are stored). Typing a
# plain numeric-looking string in directly gets auto-coerced to a real
# number, so instead compute it with TEXT() and paste back as a value -
# this keeps the result as text without touching the cell's style.
$ws.Range("D5").Formula = '=TEXT(0.81,"0.00")'
$ws.Range("D5").Copy()
$ws.Range("D5").PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteValues)

$ws.Range("F5").Formula = '=TEXT(0.977,"0.000")'
$ws.Range("F5").Copy()
$ws.Range("F5").PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteValues)

# --- Row 6 ("Nicht simultan" / "PTB"): fill in new statistics ---
$ws.Range("C6").Value2 = 0
$ws.Range("D6").Value2 = 0
$ws.Range("E6").Value2 = 0.93899999999999995
$ws.Range("F6").Value2 = 0.70099999999999996

# --- Row 7 ("Nicht simultan" / "INV"): fill in new statistics ---
$ws.Range("C7").Value2 = 0.8
$ws.Range("D7").Value2 = 0.14000000000000001
$ws.Range("E7").Value2 = 0.81899999999999995
$ws.Range("F7").Value2 = 0.45800000000000002

# Clear the clipboard marching-ants state left behind by Copy()
$excel.CutCopyMode = $false

# --- Move the active selection (matches the saved view state) ---
$ws.Range("E13").Select()
